$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Tabela1" currently spans A1:J110. Two new daily data rows were
# appended by the data bot: row 111 (2020-07-01) and row 112 (2020-07-02).

$lo = $ws.ListObjects.Item(1)

# Grow the table by two rows; this keeps the table ref / autoFilter ref and
# the sheet dimension in sync with the new data.
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Row 111: match the look of the row right above it (row 110).
$ws.Range("A110:J110").Copy() | Out-Null
$ws.Range("A111:J111").PasteSpecial(-4122) | Out-Null

# Row 112: match the look of the other "unbanded" data rows (e.g. row 100).
$ws.Range("A100:J100").Copy() | Out-Null
$ws.Range("A112:J112").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Row 111 values ---
$ws.Range("A111").Value = 44011
$ws.Range("B111").Value = 100330
$ws.Range("C111").Value = 1085
$ws.Range("D111").Value = 1600
$ws.Range("E111").Value = 15
$ws.Range("F111").Value = 8
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = "111*"
$ws.Range("J111").Value = 0

# --- Row 112 values ---
$ws.Range("A112").Value = 44012
$ws.Range("B112").Value = 101729
$ws.Range("C112").Value = 1399
$ws.Range("D112").Value = 1613
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = 8
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = "111*"
$ws.Range("J112").Value = 0

# Match the end-of-file viewport / selection state recorded in the edit.
$ws.Range("J112").Select() | Out-Null
